$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.280.14"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "2.492.10"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.95"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.58"
$ws.Range("E6").Value = "  +3.47%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("E8").Value = "  -0.34%  "

$ws.Range("D9").Value = "2.493.00"
$ws.Range("E9").Value = "  +0.87%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  +1.06%  "

$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.332"
$ws.Range("E13").Value = "  -1.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.50"
$ws.Range("E14").Value = "  -1.16%  "

$ws.Range("D15").Value = "2.919.62"

$ws.Range("D16").Value = "67.206.09"
$ws.Range("E16").Value = "  +0.81%  "

$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").Value = "2.485.58"
$ws.Range("E18").Value = "  +0.78%  "

$ws.Range("E19").Value = "  -4.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.43"
$ws.Range("E20").Value = "  -4.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.13"
$ws.Range("E21").Value = "  -2.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.03"
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.01"
$ws.Range("E23").Value = "  +0.60%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.73"
$ws.Range("E24").Value = "  -2.87%  "

$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.24"
$ws.Range("E25").Value = "  -4.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.79"
$ws.Range("E26").Value = "  -1.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").Value = "2.619.73"
$ws.Range("E29").Value = "  +0.64%  "

$ws.Range("E30").Value = "  -2.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "511.28"
$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.83"
$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("E34").Value = "  -3.21%  "

$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.32"
$ws.Range("E36").Value = "  +1.26%  "

$ws.Range("E37").Value = "  -6.34%  "

$ws.Range("E38").Value = "  +0.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.27"
$ws.Range("E39").Value = "  -3.13%  "

$ws.Range("E40").Value = "  -5.30%  "

$ws.Range("E41").Value = "  -2.65%  "

$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.84"
$ws.Range("E43").Value = "  -1.76%  "

$ws.Range("E44").Value = "  -1.27%  "

$ws.Range("E45").Value = "  -3.16%  "

$ws.Range("E46").Value = "  -1.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.07"
$ws.Range("E47").Value = "  +0.70%  "

$ws.Range("E48").Value = "  -3.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.515"
$ws.Range("E49").Value = "  -3.92%  "

$ws.Range("D50").Value = "0.0₆0251"
$ws.Range("E50").Value = "  -5.57%  "

$ws.Range("E51").Value = "  -0.58%  "
